$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'272.03"
$ws.Range("D3").Value = "'23.18"
$ws.Range("D4").Value = "'6.344"
$ws.Range("D5").Value = "'0.06340"
$ws.Range("D7").Value = "'6.783"
$ws.Range("D8").Value = "'1.401"
$ws.Range("D9").Value = "'0.8385"
$ws.Range("D10").Value = "'0.1629"
$ws.Range("D11").Value = "'0.08393"
$ws.Range("D13").Value = "'0.03165"
$ws.Range("D14").Value = "'0.09296"
$ws.Range("D15").Value = "'3.897"
$ws.Range("D16").Value = "'0.001718"
$ws.Range("D17").Value = "'0.04869"
$ws.Range("D18").Value = "'0.006316"
$ws.Range("D19").Value = "'0.005504"
$ws.Range("D20").Value = "'0.001086"
$ws.Range("D21").Value = "'0.0001495"
$ws.Range("D22").Value = "'3.738"
$ws.Range("D23").Value = "'2.352"
$ws.Range("D25").Value = "'0.3351"
$ws.Range("D26").Value = "'0.1249"
$ws.Range("D27").Value = "'0.0002674"
$ws.Range("D40").Value = "'0.04694"
$ws.Range("D41").Value = "'0.006911"
$ws.Range("D42").Value = "'0.1180"
$ws.Range("D43").Value = "'0.003449"
$ws.Range("D44").Value = "'0.01253"
$ws.Range("D45").Value = "'0.00006246"
$ws.Range("D46").Value = "'0.00000000748"
$ws.Range("D47").Value = "'0.6979"
$ws.Range("D48").Value = "'0.1252"
$ws.Range("D49").Value = "'0.00002094"
$ws.Range("D50").Value = "'0.01236"

# Reset style to Normal to clear the quote-prefix flag introduced above,
# so the cells keep matching their original (unstyled) appearance.
$ws.Range("D2:D50").Style = "Normal"
